$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 136.81682590946309
$ws.Range("C2").Value = 179.30503971052411
$ws.Range("D2").Value = 135.07188180764757
$ws.Range("E2").Value = 178.11944824163936

$ws.Range("B3").Value = 129.22963579216741
$ws.Range("C3").Value = 175.87232263042642
$ws.Range("D3").Value = 130.81180542225752
$ws.Range("E3").Value = 173.65750139831928

$ws.Range("B1:E3").Select()
